# Sync Forms export and update outputs
# This script mirrors a new Microsoft Forms response being appended to the
# "Responses" table on the worksheet, plus a small data correction on an
# existing row (N9 date got updated), a couple of UI/view tweaks, and the
# resulting widening of column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Grow the "Responses" table by one row (table ref + autofilter follow
#    automatically; dimension / table range both expand to A1:N11).
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$null = $lo.ListRows.Add()

# ---------------------------------------------------------------------
# 2) Populate the new row (row 11) with the new Forms response.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = 10

# Copy formats only from C3:D3 (style used by the e-mail columns) so the
# existing cell style (font size 12) is reused instead of creating a new one.
$ws.Range("C3:D3").Copy()
$ws.Range("C11:D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C11").Value = "achille.desbrieres@uzh.ch"
$ws.Range("D11").Value = "achille.desbrieres@uzh.ch"
$ws.Range("E11").Value = "Gizem`tTopsakal"
$ws.Range("F11").Value = "Excellent 5"
$ws.Range("G11").Value = "Excellent 5"
$ws.Range("H11").Value = "Excellent 5"
$ws.Range("I11").Value = "Excellent 5"
$ws.Range("J11").Value = "Excellent 5"
$ws.Range("K11").Value = "Excellent 5"
$ws.Range("L11").Value = "Excellent 5"
$ws.Range("M11").Value = "testestest"
$ws.Range("N11").Value = "2026-03-22T19:45:42.8842352Z"

# Match the row height Excel auto-applies to rows using the larger (size 12)
# font, as seen on the other data rows.
$ws.Rows.Item(11).RowHeight = 15.75

# ---------------------------------------------------------------------
# 3) Small data correction to an existing response (ReceivedAtUTC moved
#    from February to April).
# ---------------------------------------------------------------------
$ws.Range("N9").Value = "2026-04-22T19:44:21.8423287Z"

# ---------------------------------------------------------------------
# 4) Column B widened slightly to better fit its contents.
# ---------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 15.14

# ---------------------------------------------------------------------
# 5) Update the active selection to reflect where editing ended.
# ---------------------------------------------------------------------
$null = $ws.Range("N10").Select()
